$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Register")
$ws2 = $wb.Worksheets.Item("NewUserRegister")

# --- Sheet "Register" (sheet1), row 2 ---
# Order matters: shared-string table indices are assigned in first-seen
# write order, so write D2 (phone) before A2 (first name) before C2 (email)
# to reproduce the target's shared-string ordering.
$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "8876543210"
$ws1.Range("A2").Value = "Rgghav"
$ws1.Range("C2").Value = "Raghav11@example.com"
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:Raghav11@example.com")
$ws1.Range("C2").Style = "Hyperlink"

# Move sheet1's own cursor to C2 without leaving it as the active sheet.
$ws1.Range("C2").Select()

# --- Sheet "NewUserRegister" (sheet2), row 2 ---
$ws2.Range("B2").Value = "otte"
$ws2.Range("C2").Value = "vaibhavotte511@example.com"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "7476543215"

# Restore NewUserRegister as the active sheet/tab (it was active originally).
$ws2.Activate()
